# Week4 cancer_fraction deck: add a new "Testing" section slide that records
# the path to the sample list used for TF testing (Generate_CMD / SGE note).

$p = $ppt.ActivePresentation

# New slide re-uses the same "Title and Content" layout as the preceding
# slide (slide 13, "Predicted TF vs MAF").
$lastSlide = $p.Slides.Item($p.Slides.Count)
$layout = $lastSlide.CustomLayout
$s = $p.Slides.AddSlide($p.Slides.Count + 1, $layout)

# Leave the Title placeholder empty (no text entered on this slide).

# Content placeholder: the on-disk path to the testing sample list. Typed as
# separate chunks so the run boundaries match the way PowerPoint splits text
# around the proofing/spell-check breaks ("ghds", "processed_samples").
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.Text = "/"
$run = $body.InsertAfter("ghds")
$run = $run.InsertAfter("/groups/bioinformatics/02_DEVELOPMENT/200830_GHCNVWG_PIPELINE/results/210608_CLINICAL_SAMPLES_V0.5/")
$run = $run.InsertAfter("processed_samples")
$run = $run.InsertAfter("/List_of_samples_4TF_testing.210708.tsv")

# File this new slide under its own "Testing" section (after the existing
# "Default Section" / "Feature_reduction" sections).
$sections = $p.SectionProperties
$newSectionIndex = $sections.AddSection($sections.Count + 1, "Testing")
$s.MoveToSectionStart($newSectionIndex)
